# Pregates compensation controls in a smarter way before constructing the
# spillover matrix. Manually updated the "Comp controls" sheet so marker
# names / compensation-control FCS file names carry the "-A" (area) channel
# suffix that the parser expects, Live/Dead green FITC is renamed, and the
# PE-Cy7 / APC-H7 markers are renumbered to APC-Cy7 to match the new panel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comp controls")

# --- Marker (column B) changes, in edit order ---------------------------
$ws.Range("B13").Value = "CD127:APC-A"
$ws.Range("B4").Value  = "Live/Dead:FITC-A"
$ws.Range("B7").Value  = "CD45RA:PE-Cy7-A"
$ws.Range("B8").Value  = "CD194:PE-Cy7-A"
$ws.Range("B9").Value  = "CD27:PE-Cy7-A"
$ws.Range("B10").Value = "CD11c:PE-Cy7-A"
$ws.Range("B11").Value = "CD196:PE-Cy7-A"
$ws.Range("B14").Value = "CD8:APC-Cy7-A"
$ws.Range("B15").Value = "CD45RO:APC-Cy7-A"
$ws.Range("B17").Value = "CD3+19+20:APC-Cy7-A"

# --- Compensation control FCS file names (column A) ----------------------
$ws.Range("A4").Value  = "Compensation Controls_Live,2f,Dead FITC-A Stained Control.fcs"
$ws.Range("A5").Value  = "Compensation Controls_CD197 PE-A Stained Control.fcs"
$ws.Range("A6").Value  = "Compensation Controls_CD4 PerCP-Cy5-5-A Stained Control.fcs"
$ws.Range("A8").Value  = "Compensation Controls_CD194 PE-Cy7-A Stained Control.fcs"
$ws.Range("A12").Value = "Compensation Controls_CD38 APC-A Stained Control.fcs"
$ws.Range("A16").Value = "Compensation Controls_CD20 APC-Cy7-A Stained Control.fcs"
$ws.Range("A18").Value = "Compensation Controls_CD3 Pacific Blue-A Stained Control.fcs"
$ws.Range("A19").Value = "Compensation Controls_HLA DR Am Cyan-A Stained Control.fcs"

# --- Column B width: now auto-sized (bestFit) to fit the new, longer text --
$ws.Columns.Item(2).ColumnWidth = 19.83203125

# --- Selection moves from B20 to A3 --------------------------------------
$ws.Range("A3").Select() | Out-Null
